# Sync a new evaluation row (row 5) into the "avaliacoes_garantia" sheet,
# mirroring the existing rows' layout: stars | comment | createdAt | task.id

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A5: stars
$ws.Cells.Item(5, 1).Value = 5

# B5: comment - this record has no comment text, just like B2/B4.
# Assigning "" directly clears a cell entirely in this COM model, so we use
# the classic Excel "lone apostrophe" idiom to force an empty *text* cell,
# then reset the cell's style back to Normal so it carries no formatting
# (matching the unstyled empty cells already used for B2/B4).
$ws.Cells.Item(5, 2).Value = "'"
$ws.Cells.Item(5, 2).Style = "Normal"

# C5: createdAt (serial date/time) - reuse the same date/time display
# format already applied to the other createdAt cells in column C.
$ws.Cells.Item(5, 3).Value = 45895.62649680555
$ws.Cells.Item(5, 3).NumberFormat = "YYYY-MM-DD HH:MM:SS"

# D5: task.id (opaque id token, stored verbatim as text)
$ws.Cells.Item(5, 4).Value = "NjMzNDFkZTQtNWEwYi00MGVhLWE0YTMtZDEzMmM1YjFjNmUyOjU3MDE2"
